# Refactered the playwright into BDD cucmber pattern
#
# Update the LoginData test sheet: the "student" login value in A2 is
# replaced with a new "studentexcel" value (adds a new shared string),
# and the last active selection is moved from D11 to A10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# A2 held the shared string "student" (same value as A3); change it to a
# new distinct value "studentexcel" so it no longer shares the "student"
# string with A3.
$ws.Range("A2").Value = "studentexcel"

# Move the sheet's last selection/active cell to A10.
$ws.Range("A10").Select()
